# Insert a new data row at row 217 (a new Cilantro price observation at
# "Feria Lagunitas de Puerto Montt"), shifting every existing row from 217
# down to 304 one row down (so the old row 217 becomes row 218, ..., the old
# row 304 becomes row 305). The sheet dimension grows from A1:R304 to
# A1:R305.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 217..304 down by inserting a blank row at 217.
$ws.Rows("217").Insert(-4121, 0)

# Fill the newly inserted row 217 with the new record.
$ws.Cells.Item(217, 1).Value = 4
$ws.Cells.Item(217, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(217, 3).Value = "Los Lagos"
$ws.Cells.Item(217, 4).Value = 44726
$ws.Cells.Item(217, 5).Value = 10
$ws.Cells.Item(217, 6).Value = 100112040
$ws.Cells.Item(217, 7).Value = "Cilantro"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 180
$ws.Cells.Item(217, 11).Value = 11000
$ws.Cells.Item(217, 12).Value = 11000
$ws.Cells.Item(217, 13).Value = 11000
$ws.Cells.Item(217, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(217, 15).Value = "Región Metropolitana"
$ws.Cells.Item(217, 16).Value = 306
$ws.Cells.Item(217, 17).Value = 36
$ws.Cells.Item(217, 18).Value = "Hortaliza"
